$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unmerge the two merged ranges that are being changed so we can edit freely.
$ws.Range("C71:C77").UnMerge()
$ws.Range("D71:D77").UnMerge()
$ws.Range("C78:C82").UnMerge()
$ws.Range("D78:D82").UnMerge()

# Row 78 becomes a new "3.5.8 Navigation and communication" detail row
# (was the "3.6.0 Transport and communication" header row).
$ws.Range("C78").Value = $null
$ws.Range("D78").Value = $null
$ws.Range("E78").Value = "3.5.8"
$ws.Range("F78").Value = "Navigation and communication"

# Row 79 becomes the new "3.6.0 Transport" header row
# (renamed from "Transport and communication"), with sub item 3.6.1 unchanged.
$ws.Range("C79").Value = "3.6.0"
$ws.Range("D79").Value = "Transport"
$ws.Range("E79").Value = "3.6.1"
$ws.Range("F79").Value = "Airports/aerodromes"

# Rows 80-82 shift their E numbering down by one (F labels unchanged text-wise,
# but shift up one row as a consequence of the renumbering below).
$ws.Range("E80").Value = "3.6.2"
$ws.Range("F80").Value = "Roads"

$ws.Range("E81").Value = "3.6.3"
$ws.Range("F81").Value = "Railways"

$ws.Range("E82").Value = "3.6.4"
$ws.Range("F82").Value = "Ports and water transport"

# Re-merge the adjusted ranges: Utilities group now spans 71:78, and the
# Transport header/group now spans 79:82.
$ws.Range("C71:C78").Merge()
$ws.Range("D71:D78").Merge()
$ws.Range("C79:C82").Merge()
$ws.Range("D79:D82").Merge()
